# Update checklist criteria text (remove leading "Category | " prefix;
# row 14 also changes "roteou" -> "transferiu") on both the Checklist
# and Config sheets, which mirror the same criteria list in column B.
$wb = $excel.ActiveWorkbook
$wsChecklist = $wb.Worksheets.Item("Checklist")
$wsConfig = $wb.Worksheets.Item("Config")

$criteria = @(
    'Não cumprimentou corretamente (padrão regional/linguístico)',
    'Não reconheceu todas as preocupações do cliente',
    'Não coletou informações relevantes (comentários/ferramentas)',
    'Não fez perguntas relevantes',
    'Resolução fora da Base de Conhecimento (KB)',
    'Informações imprecisas, irrelevantes ou incompletas',
    'Não abordou todas as preocupações do cliente',
    'Não resolveu questões adicionais da conta',
    'Ações externas (voltadas ao cliente)',
    'Ações internas (processos internos)',
    'Não realizou ajustes monetários necessários (dentro do limite)',
    'Não transferiu corretamente o contato',
    'Não fechou contatos duplicados',
    'Não compartilhou detalhes internos conforme diretrizes',
    'Não escalou corretamente conforme KB',
    'Não compartilhou detalhes da escalação corretamente',
    'Tipo de problema selecionado incorretamente',
    'Não seguiu diretrizes de tempo de espera',
    'Não definiu expectativas corretamente',
    'Não orientou nem educou o cliente',
    'Não verificou necessidades adicionais',
    'Falta de empatia no atendimento',
    'Não personalizou a situação do cliente',
    'Erros de gramática, ortografia ou formatação',
    'Uso excessivo de respostas salvas / tom robótico',
    'Interrompeu, cortou a fala ou apressou o cliente',
    'Comprometeu a privacidade do cliente (Segurança/PCI)',
    'Violou confidencialidade da empresa (informações internas)',
    'Comunicação inadequada (comentários negativos, rudeza) ',
    'Não escalou corretamente questão de jogo responsável',
    'Não escalou corretamente ameaças legais/regulatórias'
)

for ($i = 0; $i -lt $criteria.Count; $i++) {
    $row = 3 + $i
    $wsChecklist.Range("B$row").Value = $criteria[$i]
    $wsConfig.Range("B$row").Value = $criteria[$i]
}

# Update the saved selection/scroll state on each sheet to match the
# author's last-known cursor position when the file was re-saved.
$wsChecklist.Range("C30").Select()

$wsConfig.Activate()
$wsConfig.Range("B3").Select()
